$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows with name/link swaps (34<->35, 45<->46) ---
$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D34").Value = "'0.07845"
$ws.Range("E34").Value = "  +1.30%  "

$ws.Range("B35").Value = "WEMIXTOKEN"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "'1.575"
$ws.Range("E35").Value = "  +22.94%  "

$ws.Range("B45").Value = "PancakeSwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D45").Value = "'3.591"
$ws.Range("E45").Value = "  +1.47%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'12.46"
$ws.Range("E46").Value = "  +1.44%  "

# --- Price / Volume updates ---
$ws.Range("D2").Value = "'20.544.13"
$ws.Range("E2").Value = "  +2.80%  "
$ws.Range("D3").Value = "'1.467.68"
$ws.Range("E3").Value = "  +3.44%  "
$ws.Range("E4").Value = "  +0.66%  "
$ws.Range("D5").Value = "'0.9356"
$ws.Range("E5").Value = "  -6.62%  "
$ws.Range("D6").Value = "'280.99"
$ws.Range("E6").Value = "  +2.45%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.3197"
$ws.Range("E8").Value = "  +3.80%  "
$ws.Range("D9").Value = "'41.31"
$ws.Range("E9").Value = "  +4.52%  "
$ws.Range("E10").Value = "  +5.35%  "
$ws.Range("D11").Value = "'0.06686"
$ws.Range("E11").Value = "  +1.36%  "
$ws.Range("D12").Value = "'1.004"
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("D13").Value = "'5.592"
$ws.Range("E13").Value = "  +3.17%  "
$ws.Range("D14").Value = "'18.29"
$ws.Range("E14").Value = "  +6.85%  "
$ws.Range("D15").Value = "'6.246"
$ws.Range("E15").Value = "  +1.13%  "
$ws.Range("D16").Value = "'1.476.10"
$ws.Range("E16").Value = "  +3.67%  "
$ws.Range("D17").Value = "'0.00001038"
$ws.Range("E17").Value = "  +2.97%  "
$ws.Range("D18").Value = "'0.9384"
$ws.Range("E18").Value = "  -6.29%  "
$ws.Range("D19").Value = "'0.05739"
$ws.Range("E19").Value = "  -0.55%  "
$ws.Range("D20").Value = "'72.28"
$ws.Range("E20").Value = "  -2.81%  "
$ws.Range("D21").Value = "'5.700"
$ws.Range("E21").Value = "  +1.30%  "
$ws.Range("D22").Value = "'14.82"
$ws.Range("E22").Value = "  +2.12%  "
$ws.Range("D23").Value = "'11.23"
$ws.Range("E23").Value = "  +2.11%  "
$ws.Range("D24").Value = "'2.287"
$ws.Range("E24").Value = "  -2.05%  "
$ws.Range("D25").Value = "'20.664.03"
$ws.Range("E25").Value = "  +3.34%  "
$ws.Range("D26").Value = "'2.306"
$ws.Range("E26").Value = "  +0.75%  "
$ws.Range("D27").Value = "'138.19"
$ws.Range("E27").Value = "  -0.78%  "
$ws.Range("D28").Value = "'17.62"
$ws.Range("E28").Value = "  +3.78%  "
$ws.Range("D29").Value = "'1.637.60"
$ws.Range("E29").Value = "  +3.51%  "
$ws.Range("D30").Value = "'113.89"
$ws.Range("E30").Value = "  +4.29%  "
$ws.Range("D31").Value = "'3.953"
$ws.Range("E31").Value = "  +3.24%  "
$ws.Range("D32").Value = "'5.296"
$ws.Range("E32").Value = "  -1.96%  "
$ws.Range("D33").Value = "'0.8546"
$ws.Range("E33").Value = "  -3.18%  "
$ws.Range("D36").Value = "'0.06110"
$ws.Range("E36").Value = "  +6.29%  "
$ws.Range("D37").Value = "'4.938"
$ws.Range("E37").Value = "  +2.96%  "
$ws.Range("D38").Value = "'10.72"
$ws.Range("E38").Value = "  -1.95%  "
$ws.Range("D39").Value = "'0.02071"
$ws.Range("E39").Value = "  +1.05%  "
$ws.Range("D40").Value = "'1.124"
$ws.Range("E40").Value = "  +4.90%  "
$ws.Range("D41").Value = "'0.1907"
$ws.Range("E41").Value = "  -1.46%  "
$ws.Range("D42").Value = "'0.9525"
$ws.Range("E42").Value = "  -4.84%  "
$ws.Range("D43").Value = "'7.505"
$ws.Range("E43").Value = "  -11.36%  "
$ws.Range("D44").Value = "'0.5414"
$ws.Range("E44").Value = "  +1.75%  "
$ws.Range("D47").Value = "'122.23"
$ws.Range("E47").Value = "  +11.66%  "
$ws.Range("D48").Value = "'0.5336"
$ws.Range("E48").Value = "  +3.91%  "
$ws.Range("D49").Value = "'1.830"
$ws.Range("E49").Value = "  +1.55%  "
$ws.Range("D50").Value = "'0.06462"
$ws.Range("E50").Value = "  +4.97%  "
$ws.Range("D51").Value = "'1.047"
$ws.Range("E51").Value = "  -0.20%  "
